# Update odds/results data per re-scrape on 06-11-2023 02:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 282: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F282').Value = 'Athletico-PR'
$ws.Range('G282').Value = 3
$ws.Range('H282').Value = 'America MG'
$ws.Range('I282').Value = 2
$ws.Range('J282').Value = 1.58
$ws.Range('K282').Value = '22/10/2023 22:42'
$ws.Range('L282').Value = 1.58
$ws.Range('M282').Value = '25/10/2023 23:27'
$ws.Range('N282').Value = 4.26
$ws.Range('O282').Value = '22/10/2023 22:42'
$ws.Range('P282').Value = 4.31
$ws.Range('Q282').Value = '25/10/2023 23:33'
$ws.Range('R282').Value = 6.01
$ws.Range('S282').Value = '22/10/2023 22:42'
$ws.Range('T282').Value = 6.03
$ws.Range('U282').Value = '25/10/2023 23:30'
$ws.Range('V282').Value = 'https://www.betexplorer.com/football/brazil/serie-a/athletico-pr-america-mg/AXnrQVpd/'

# Row 283: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F283').Value = 'Bragantino'
$ws.Range('G283').Value = 1
$ws.Range('H283').Value = 'Atletico-MG'
$ws.Range('I283').Value = 2
$ws.Range('J283').Value = 1.96
$ws.Range('K283').Value = '22/10/2023 22:42'
$ws.Range('L283').Value = 1.88
$ws.Range('M283').Value = '25/10/2023 23:59'
$ws.Range('N283').Value = 3.5
$ws.Range('O283').Value = '22/10/2023 22:42'
$ws.Range('P283').Value = 3.47
$ws.Range('Q283').Value = '25/10/2023 23:59'
$ws.Range('R283').Value = 4.2
$ws.Range('S283').Value = '22/10/2023 22:42'
$ws.Range('T283').Value = 4.8
$ws.Range('U283').Value = '25/10/2023 23:59'
$ws.Range('V283').Value = 'https://www.betexplorer.com/football/brazil/serie-a/bragantino-atletico-mg/Cd77WX7S/'

# Row 284: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F284').Value = 'Fluminense'
$ws.Range('G284').Value = 5
$ws.Range('H284').Value = 'Goias'
$ws.Range('I284').Value = 3
$ws.Range('J284').Value = 1.53
$ws.Range('K284').Value = '22/10/2023 22:42'
$ws.Range('L284').Value = 1.61
$ws.Range('M284').Value = '25/10/2023 23:34'
$ws.Range('N284').Value = 4.28
$ws.Range('O284').Value = '22/10/2023 22:42'
$ws.Range('P284').Value = 4.15
$ws.Range('Q284').Value = '25/10/2023 23:34'
$ws.Range('R284').Value = 6.81
$ws.Range('S284').Value = '22/10/2023 22:42'
$ws.Range('T284').Value = 5.9
$ws.Range('U284').Value = '25/10/2023 23:34'
$ws.Range('V284').Value = 'https://www.betexplorer.com/football/brazil/serie-a/fluminense-goias/rHFjzG83/'

# Row 287: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F287').Value = 'Gremio'
$ws.Range('G287').Value = 3
$ws.Range('H287').Value = 'Flamengo RJ'
$ws.Range('I287').Value = 2
$ws.Range('J287').Value = 2.54
$ws.Range('K287').Value = '22/10/2023 20:12'
$ws.Range('L287').Value = 4.14
$ws.Range('M287').Value = '26/10/2023 02:27'
$ws.Range('N287').Value = 3.37
$ws.Range('O287').Value = '22/10/2023 20:12'
$ws.Range('P287').Value = 3.4
$ws.Range('Q287').Value = '26/10/2023 02:22'
$ws.Range('R287').Value = 2.93
$ws.Range('S287').Value = '22/10/2023 20:12'
$ws.Range('T287').Value = 2.03
$ws.Range('U287').Value = '26/10/2023 02:27'
$ws.Range('V287').Value = 'https://www.betexplorer.com/football/brazil/serie-a/gremio-flamengo-rj/WtlvRBVk/'

# Row 288: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F288').Value = 'Cuiaba'
$ws.Range('G288').Value = 0
$ws.Range('H288').Value = 'Corinthians'
$ws.Range('I288').Value = 1
$ws.Range('J288').Value = 2.03
$ws.Range('K288').Value = '22/10/2023 22:42'
$ws.Range('L288').Value = 2.23
$ws.Range('M288').Value = '26/10/2023 02:29'
$ws.Range('N288').Value = 3.2
$ws.Range('O288').Value = '22/10/2023 22:42'
$ws.Range('P288').Value = 3.1
$ws.Range('Q288').Value = '26/10/2023 02:27'
$ws.Range('R288').Value = 4.27
$ws.Range('S288').Value = '22/10/2023 22:42'
$ws.Range('T288').Value = 3.88
$ws.Range('U288').Value = '26/10/2023 02:29'
$ws.Range('V288').Value = 'https://www.betexplorer.com/football/brazil/serie-a/cuiaba-corinthians/MLgTmZx3/'

# Row 291: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F291').Value = 'America MG'
$ws.Range('G291').Value = 3
$ws.Range('H291').Value = 'Gremio'
$ws.Range('I291').Value = 4
$ws.Range('J291').Value = 2.54
$ws.Range('K291').Value = '26/10/2023 02:42'
$ws.Range('L291').Value = 2.32
$ws.Range('M291').Value = '28/10/2023 23:59'
$ws.Range('N291').Value = 3.43
$ws.Range('O291').Value = '26/10/2023 02:42'
$ws.Range('P291').Value = 3.56
$ws.Range('Q291').Value = '28/10/2023 23:59'
$ws.Range('R291').Value = 2.89
$ws.Range('S291').Value = '26/10/2023 02:42'
$ws.Range('T291').Value = 3.16
$ws.Range('U291').Value = '28/10/2023 23:59'
$ws.Range('V291').Value = 'https://www.betexplorer.com/football/brazil/serie-a/america-mg-gremio/MsV30CM8/'

# Row 292: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F292').Value = 'Palmeiras'
$ws.Range('G292').Value = 1
$ws.Range('H292').Value = 'Bahia'
$ws.Range('I292').Value = 0
$ws.Range('J292').Value = 1.51
$ws.Range('K292').Value = '26/10/2023 01:12'
$ws.Range('L292').Value = 1.34
$ws.Range('M292').Value = '28/10/2023 23:52'
$ws.Range('N292').Value = 4.44
$ws.Range('O292').Value = '26/10/2023 01:12'
$ws.Range('P292').Value = 5.16
$ws.Range('Q292').Value = '28/10/2023 23:53'
$ws.Range('R292').Value = 6.84
$ws.Range('S292').Value = '26/10/2023 01:12'
$ws.Range('T292').Value = 10.94
$ws.Range('U292').Value = '28/10/2023 23:53'
$ws.Range('V292').Value = 'https://www.betexplorer.com/football/brazil/serie-a/palmeiras-bahia/CAntoXML/'

# Row 313: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F313').Value = 'Fortaleza'
$ws.Range('G313').Value = 0
$ws.Range('H313').Value = 'Flamengo RJ'
$ws.Range('I313').Value = 2
$ws.Range('J313').Value = 2.84
$ws.Range('K313').Value = '02/11/2023 01:42'
$ws.Range('L313').Value = 2.58
$ws.Range('M313').Value = '05/11/2023 19:58'
$ws.Range('N313').Value = 3.27
$ws.Range('O313').Value = '02/11/2023 01:42'
$ws.Range('P313').Value = 3.16
$ws.Range('Q313').Value = '05/11/2023 19:52'
$ws.Range('R313').Value = 2.63
$ws.Range('S313').Value = '02/11/2023 01:42'
$ws.Range('T313').Value = 3.09
$ws.Range('U313').Value = '05/11/2023 19:58'
$ws.Range('V313').Value = 'https://www.betexplorer.com/football/brazil/serie-a/fortaleza-flamengo-rj/6BfQDMdP/'

# Row 314: swap in re-scraped match data (same date, reordered by site)
$ws.Range('F314').Value = 'Bragantino'
$ws.Range('G314').Value = 1
$ws.Range('H314').Value = 'Corinthians'
$ws.Range('I314').Value = 0
$ws.Range('J314').Value = 1.58
$ws.Range('K314').Value = '02/11/2023 22:12'
$ws.Range('L314').Value = 1.6
$ws.Range('M314').Value = '05/11/2023 19:53'
$ws.Range('N314').Value = 4.08
$ws.Range('O314').Value = '02/11/2023 22:12'
$ws.Range('P314').Value = 4.22
$ws.Range('Q314').Value = '05/11/2023 19:58'
$ws.Range('R314').Value = 6.35
$ws.Range('S314').Value = '02/11/2023 22:12'
$ws.Range('T314').Value = 5.87
$ws.Range('U314').Value = '05/11/2023 19:57'
$ws.Range('V314').Value = 'https://www.betexplorer.com/football/brazil/serie-a/bragantino-corinthians/hl4Wiuz6/'

# New row 315 (Coritiba vs Goias) - copy formatting from row 314, then set values
$ws.Range('A314:V314').Copy()
$ws.Range('A315:V315').PasteSpecial(-4122)
$ws.Range('D315').NumberFormat = '@'

$ws.Range('A315').Value = 314
$ws.Range('B315').Value = 'brazil'
$ws.Range('C315').Value = 'serie-a'
$ws.Range('D315').Value = '2023'
$ws.Range('E315').Value = 45235.9375
$ws.Range('F315').Value = 'Coritiba'
$ws.Range('G315').Value = 0
$ws.Range('H315').Value = 'Goias'
$ws.Range('I315').Value = 1
$ws.Range('J315').Value = 2.29
$ws.Range('K315').Value = '02/11/2023 22:12'
$ws.Range('L315').Value = 2.7
$ws.Range('M315').Value = '05/11/2023 22:28'
$ws.Range('N315').Value = 3.29
$ws.Range('O315').Value = '02/11/2023 22:12'
$ws.Range('P315').Value = 3.33
$ws.Range('Q315').Value = '05/11/2023 22:28'
$ws.Range('R315').Value = 3.44
$ws.Range('S315').Value = '02/11/2023 22:12'
$ws.Range('T315').Value = 2.8
$ws.Range('U315').Value = '05/11/2023 22:28'
$ws.Range('V315').Value = 'https://www.betexplorer.com/football/brazil/serie-a/coritiba-goias/46GskwKO/'
